# [Fix]: exclusion of 8 redundant metrics
#
# Remove the 8 redundant metric rows (MBRAE, UMBRAE, STDAPE, RMSPE, MRE,
# MRAE, MDRAE, GMRAE) from both the "LMN" and "STS" worksheets, then
# renumber the remaining "ID" column sequentially.

$wb = $excel.ActiveWorkbook

# Row numbers (1-based, matching the worksheet's current row numbers) of
# the metrics to drop. Identical layout on both sheets.
$rowsToDelete = @(33, 32, 31, 29, 24, 23, 13, 12)

foreach ($ws in $wb.Worksheets) {
    foreach ($r in $rowsToDelete) {
        $ws.Rows($r).Delete()
    }

    # Renumber column A (the "ID" values) sequentially for the surviving
    # data rows (row 1 is the header, data starts at row 2).
    $lastRow = $ws.UsedRange.Rows.Count
    for ($r = 2; $r -le $lastRow; $r++) {
        $ws.Cells.Item($r, 1).Value = $r - 1
    }
}
